$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP")
foreach ($col in $cols) {
  $src = $ws.Range($col + "20")
  $dst = $ws.Range($col + "21")
  $src.Copy()
  $dst.PasteSpecial(-4122)
}
Write-Host "done"
